# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handed-off files have now been handed back and are in sync with en-US.

$wb = $excel.ActiveWorkbook

# Colors / constants matching the workbook's existing "HyperLink" cell style
# (underline, font color FF6495ED == RGB(100,149,237))
$hyperlinkColor = 15570276   # RGB(0x64, 0x95, 0xED) packed as BGR long for COM
$underlineSingle = 2         # xlUnderlineStyleSingle

function Set-HyperlinkCell {
    param(
        $ws,
        [string]$cellRef,
        [string]$url,
        [string]$displayText
    )
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText)
    $ws.Range($cellRef).Font.Underline = $underlineSingle
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# 1) Status changes from "Ready for handoff" to "Handed back: in sync with
#    en-US" everywhere it is used: the Overview sheet (per file/language)
#    as well as the "Status" column on each per-language detail sheet.
# ---------------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: populate "Latest Target File" (F) and "Latest Handback
#    File" (G) with hyperlinks, and set "Latest Handback DateTime" (H).
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

Set-HyperlinkCell $zhcn "F2" "https://github.com/OpenLocalizationTest/oltest/blob/8a4f466e80f8efc45e009b0d8e89a175b5bd2386/e2e/6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.md" "6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.md"
Set-HyperlinkCell $zhcn "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/219843776f7c3ed4ac27f1a0a9b9b5c74e59a620/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.61559ace879255c3ddec3ea340c6fcf55e97712c.zh-cn.xlf" "6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.61559ace879255c3ddec3ea340c6fcf55e97712c.zh-cn.xlf"

Set-HyperlinkCell $zhcn "F3" "https://github.com/OpenLocalizationTest/oltest/blob/8a4f466e80f8efc45e009b0d8e89a175b5bd2386/e2e/9b430b49-e063-4173-a17b-23c0f65f5001.md" "9b430b49-e063-4173-a17b-23c0f65f5001.md"
Set-HyperlinkCell $zhcn "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/219843776f7c3ed4ac27f1a0a9b9b5c74e59a620/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/9b430b49-e063-4173-a17b-23c0f65f5001.03b0e92c69035084e9ef4a6a2a5af615c4ad1872.zh-cn.xlf" "9b430b49-e063-4173-a17b-23c0f65f5001.03b0e92c69035084e9ef4a6a2a5af615c4ad1872.zh-cn.xlf"

$zhcn.Range("H2").NumberFormat = "@"
$zhcn.Range("H2").Value = "2016-03-20 10:11:51"
$zhcn.Range("H3").NumberFormat = "@"
$zhcn.Range("H3").Value = "2016-03-20 10:11:51"

# ---------------------------------------------------------------------------
# 3) de-de sheet: populate "Latest Target File" (F) and "Latest Handback
#    File" (G) with hyperlinks, and set "Latest Handback DateTime" (H).
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

Set-HyperlinkCell $dede "F2" "https://github.com/OpenLocalizationTest/oltest/blob/8a4f466e80f8efc45e009b0d8e89a175b5bd2386/e2e/6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.md" "6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.md"
Set-HyperlinkCell $dede "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5385052427aa7bd014a3924c23f967754d154f42/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.61559ace879255c3ddec3ea340c6fcf55e97712c.de-de.xlf" "6ce7441b-65ca-4454-9e6f-1b40ceb02e8f.61559ace879255c3ddec3ea340c6fcf55e97712c.de-de.xlf"

Set-HyperlinkCell $dede "F3" "https://github.com/OpenLocalizationTest/oltest/blob/8a4f466e80f8efc45e009b0d8e89a175b5bd2386/e2e/9b430b49-e063-4173-a17b-23c0f65f5001.md" "9b430b49-e063-4173-a17b-23c0f65f5001.md"
Set-HyperlinkCell $dede "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5385052427aa7bd014a3924c23f967754d154f42/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/9b430b49-e063-4173-a17b-23c0f65f5001.03b0e92c69035084e9ef4a6a2a5af615c4ad1872.de-de.xlf" "9b430b49-e063-4173-a17b-23c0f65f5001.03b0e92c69035084e9ef4a6a2a5af615c4ad1872.de-de.xlf"

$dede.Range("H2").NumberFormat = "@"
$dede.Range("H2").Value = "2016-03-20 10:11:56"
$dede.Range("H3").NumberFormat = "@"
$dede.Range("H3").Value = "2016-03-20 10:11:56"
